$wb = $excel.ActiveWorkbook

# Update the "DateProd" (column B, row 2) timestamp on each sheet to reflect
# the latest Katalon test execution run (BWP Object Repository / VRelay fixes).

$ws = $wb.Worksheets.Item("AddDeleteRole")
$ws.Range("B2").Value = "Thu Sep 04 06:14:25 IST 2025"

$ws = $wb.Worksheets.Item("SearchRole")
$ws.Range("B2").Value = "Thu Sep 04 06:15:12 IST 2025"

$ws = $wb.Worksheets.Item("CreateUser")
$ws.Range("B2").Value = "Thu Sep 04 06:15:40 IST 2025"

$ws = $wb.Worksheets.Item("FindUser")
$ws.Range("B2").Value = "Thu Sep 04 06:16:22 IST 2025"

$ws = $wb.Worksheets.Item("ModifyUser")
$ws.Range("B2").Value = "Thu Sep 04 06:16:55 IST 2025"

$ws = $wb.Worksheets.Item("ModifyUserPwd")
$ws.Range("B2").Value = "Thu Sep 04 06:17:37 IST 2025"

$ws = $wb.Worksheets.Item("FindCaseUser")
$ws.Range("B2").Value = "Thu Sep 04 06:18:44 IST 2025"
